$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing trailer links (J column) ---
$ws.Cells.Item(3, 10).Value = "www.youtube.com/watch?v=bpc7ScQkIvs"
$ws.Cells.Item(9, 10).Value = "www.youtube.com/watch?v=JwG3kZSox2g"
$ws.Cells.Item(18, 10).Value = "www.youtube.com/watch?v=y1wBNlYm_YM"
$ws.Cells.Item(21, 10).Value = "www.youtube.com/watch?v=cgoTUHHGHJk"

# --- Append new rows 25-51 ---
$newRows = @(
    @(23, "15.10.2020", "Astronaut ", "Shelagh McLeod", "Richard Dreyfuss, Lyriq Bent, Krista Bridges, Colm Feore, Richie Lawrence", "https://www.kino-zeit.de/node/50228", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-07/astronaut_3.jpg", "JETS Filmverleih & Vertrieb", "Kanada", "www.youtube.com/watch?v=FLnE15vAHcA"),
    @(24, "15.10.2020", "Bruno - Die lange Heimreise ", "Karl Golden", "Diarmaid Murtagh, Woody Norman, Seun Shote, Scarlett Alice Johnson", "https://www.kino-zeit.de/node/48854", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-08/bruno_7_daniel_andizzy_searching_copyrightfilmperlen.jpg", "Filmperlen", "Großbritannien", "www.youtube.com/watch?v=_rpLj7zYEYI"),
    @(25, "18.10.2020", "Chernobyl: Abyss ", "Danila Kozlovsky", "Oksana Akinshina, Danila Kozlowskj, Filipp Avdeew", "https://www.kino-zeit.de/node/50781", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-09/chernobyl_abyss_2.jpg", "Capelight Pictures / Central", "Russland", "www.youtube.com/watch?v=amiO_cqRKis"),
    @(26, "15.10.2020", "Der Bär in mir ", "Roman Droux", $null, "https://www.kino-zeit.de/node/50292", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-07/der_bar_in_mir_2.jpg", "MFA + Filmdistribution", "Schweiz", "www.youtube.com/watch?v=JYAGKtI94eI"),
    @(27, "15.10.2020", "Der geheime Garten ", "Marc Munden", "Colin Firth, Maeve Dermody, Julie Walters, Dixie Egerickx", "https://www.kino-zeit.de/node/47584", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-02/dergeheimegarten_003.jpg", "STUDIOCANAL GmbH Filmverleih", "Großbritannien, Frankreich", "www.youtube.com/watch?v=ucZlnduCbL0"),
    @(28, "16.10.2020", "I Am Greta ", "Nathan Grossman", "Greta Thunberg", "https://www.kino-zeit.de/node/50550", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-08/Greta_Biennale.jpg", "Filmwelt Verleihagentur", "Schweden", "www.youtube.com/watch?v=xDdEWkA15Rg"),
    @(29, "17.10.2020", "Making Montgomery Clift ", "Robert Anderson Clift, Hillary Demmon", $null, "https://www.kino-zeit.de/node/47343", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2019-08/Making-Montgomery_Clift_filmbild1.jpg", "missingFILMs ", "USA", "www.youtube.com/watch?v=eJtdIcOR2bM"),
    @(30, "15.10.2020", "Martin Margiela - Mythos der Mode ", "Reiner Holzemer", $null, "https://www.kino-zeit.de/node/48816", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-02/martin_margiela_in_his_own_words_filmstill.jpg", "Filmwelt Verleihagentur", "Deutschland, Belgien", "www.youtube.com/watch?v=EvrQX4POKtg"),
    @(31, "15.10.2020", "Moskau Einfach! ", "Micha Lewinsky", "Miriam Stein, Oriana Schrage, Urs Jucker, Michael Maertens, Mike Müller", "https://www.kino-zeit.de/node/48285", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2019-12/moskau_einfach_filmbild1.jpg", "Arsenal Filmverleih GmbH ", "Schweiz", "www.youtube.com/watch?v=276LLcUGkq8"),
    @(32, "15.10.2020", "Mrs. Taylor's Singing Club ", "Peter Cattaneo", "Kristin Scott Thomas, Sharon Horgan, Emma Lowndes, Amy James-Kelly", "https://www.kino-zeit.de/node/48466", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-01/military-wives-filmbild1.jpg", "Leonine", "Großbritannien", "www.youtube.com/watch?v=BmDav7mUjMc"),
    @(33, "15.10.2020", "Oeconomia ", "Carmen Losmann", $null, "https://www.kino-zeit.de/node/48739", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-02/Oeconomia.jpg", "Neue Visionen Filmverleih GmbH", "Deutschland", "www.youtube.com/watch?v=RKT_g14MSus"),
    @(34, "15.10.2020", "Rojo - Wenn alle schweigen, ist keiner unschuldig ", "Benjámin Naishtat", "Darío Grandinetti, Andrea Frigerio, Alfredo Castro", "https://www.kino-zeit.de/node/43353", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2018-08/rojo-filmbild1.jpg", "Cine Global Filmverleih", "Argentinien, Deutschland", "www.youtube.com/watch?v=4QhqxLA3aH4"),
    @(35, "15.10.2020", "Sag du es mir ", "Michael Fetter Nathansky", "Christina Große, Marc Ben Puch, Gisa Flake, Hendrik Arnst, Patrick Heinrich, Walid Al-Atiyat", "https://www.kino-zeit.de/node/48412", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-08/sagduesmir_filmstill_04.jpg", "missingFILMs ", "Deutschland", "www.youtube.com/watch?v=UerHnKI8xZs"),
    @(36, "15.10.2020", "Tod auf dem Nil ", "Kenneth Branagh", "Kenneth Branagh, Gal Gadot, Armie Hammer, Annette Bening, Sophie Okonedo, Russell Brand, Dawn French", "https://www.kino-zeit.de/node/49615", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-08/tod-auf-dem-nil-1.jpg", "The Walt Disney Company (Germany) GmbH ", "USA", "www.youtube.com/watch?v=npTTR2LDFGI"),
    @(37, "15.10.2020", "Von Liebe und Krieg ", "Kasper Torsting", "Tom Wlaschiha, Ulrich Thomsen, Thure Lindhardt, Natalie Madueño, Sebastian Jessen", "https://www.kino-zeit.de/node/50268", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-09/l-u-k.jpg", "Tamtam Film", "Dänemark, Deutschland, Tschechische Republik", "www.youtube.com/watch?v=t406bUnfYGo"),
    @(38, "22.10.2020", "Bohnenstange ", "Kantemir Balagov", $null, "https://www.kino-zeit.de/node/46389", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2019-05/beanpole-filmbild1.jpg", "eksystent distribution filmverleih", "Russland", "www.youtube.com/watch?v=PpUx8cxC4RQ"),
    @(39, "22.10.2020", "Cortex ", "Moritz Bleibtreu", "Moritz Bleibtreu, Jannis Niewöhner, Wotan Wilke Möhring, Nadja Uhl, Nicholas Ofczarek", "https://www.kino-zeit.de/node/50644", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-08/rev-1-ger_cortex_szenenbilder_dscf1872_high_res_jpeg.jpeg", "Warner Bros.", "Deutschland", "www.youtube.com/watch?v=1rPsZmICIn4"),
    @(40, "22.10.2020", "Die Stimme des Regenwaldes ", "Niklaus Hilber", "Sven Schelker, Charlotte Heinimann, Matthew Crowley, Daniel Ludwig", "https://www.kino-zeit.de/node/47764", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2019-10/bruno-manser-filmbild1.jpg", "Camino Filmverleih GmbH", "Schweiz, Österreich", "www.youtube.com/watch?v=8B8O2qNezmE"),
    @(41, "22.10.2020", "Ema ", "Pablo Larrain", "Gael Garcia Bernal, Mariana Di Girolamo, Santiago Cabrera, Mariana Loyola", "https://www.kino-zeit.de/node/44753", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-08/ema_filmstill_09.jpg", "Koch Films GmbH, Studiocanal", "Chile", "www.youtube.com/watch?v=Ft79LOIXlfo"),
    @(42, "22.10.2020", "Kajillionaire ", "Miranda July", "Evan Rachel Wood, Gina Rodriguez, Richard Jenkins, Debra Winger", "https://www.kino-zeit.de/node/48596", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-01/Kajillionaire-filmstill.jpg", "Universal Pictures International Germany GmbH", "USA", "www.youtube.com/watch?v=xiMPCevu8Wk"),
    @(43, "22.10.2020", "Komm und sieh ", "Elem Klimow", $null, "https://www.kino-zeit.de/node/50356", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-07/Komm-und-sieh-bild1.jpg", "Drop-Out Cinema eG", "UdSSR", "www.youtube.com/watch?v=8Ab60vCIhh0"),
    @(44, "22.10.2020", "Mein Liebhaber, der Esel & ich ", "Caroline Vignal", "Laure Calamy, Benjamin Lavernhe, Olivia Côte, Marc Fraize", "https://www.kino-zeit.de/node/50818", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-09/meinliebhaberdereselundich-1.jpg", "Capelight Pictures / Wild Bunch", "Frankreich", "www.youtube.com/watch?v=s5hgoO27log"),
    @(45, "22.10.2020", "Regeln am Band, bei hoher Geschwindigkeit ", "Yulia Lokshina", $null, "https://www.kino-zeit.de/node/48463", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-01/REGELNAMBAND_Still_3_%C2%AEwirFILM.png", "jip Film & Verleih", "Deutschland", "www.youtube.com/watch?v=nMr_Fwe2yus"),
    @(46, "22.10.2020", "The Beach House ", "Jeffrey A. Brown", "Liana Liberato, Noah Le Gros, Jake Weber, Maryann Nagel", "https://www.kino-zeit.de/node/50175", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-10/thebeachhouse_filmstill_02.jpg", "Koch Films GmbH, 24 Bilder", "USA", "www.youtube.com/watch?v=1SefWBdRjDU"),
    @(47, "24.10.2020", "The Great Green Wall ", "Jared P. Scott", $null, "https://www.kino-zeit.de/node/48568", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-01/the-great-green-wall-filmstill5.jpg", "Weltkino Filmverleih GmbH", "Großbritannien", "www.youtube.com/watch?v=66EUzT2bG4c"),
    @(48, "22.10.2020", "The Mortuary - Jeder Tod hat eine Geschichte ", "Ryan Spindell", "Clancy Brown, Caitlin Fisher, Christine Kilmer, Ema Horvath", "https://www.kino-zeit.de/node/50782", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-10/mortuary_01.jpg", "Capelight Pictures, Central Filmverleih", "USA", "www.youtube.com/watch?v=-XZ6Mmbo9SA"),
    @(49, "22.10.2020", "Winterreise ", "Anders Østergaard, Erzsébet Rácz", "Bruno Ganz", "https://www.kino-zeit.de/node/49515", "https://www.kino-zeit.de/sites/default/files/styles/2560_x_1440_header_/public/2020-04/winterreise1.jpg", "Real Fiction Filmverleih", "Dänemark, Deutschland", "www.youtube.com/watch?v=TTU6-nQp5oI")
)

$startRow = 25
for ($idx = 0; $idx -lt $newRows.Count; $idx++) {
    $r = $startRow + $idx
    $row = $newRows[$idx]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    if ($row[4] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $row[4]
    }
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
}

# --- Apply the A-column border/bold/center-top style used by existing rows ---
$ws.Range("A24").Copy()
$ws.Range("A25:A51").PasteSpecial(-4122)